$wb = $excel.ActiveWorkbook

# Add new worksheet at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "BlankRow"
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# Populate header row
$newSheet.Range("A1").Value = "Id"
$newSheet.Range("B1").Value = "Name"

# Data rows with blank separator rows
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "Hoge"

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "Fuga"

$newSheet.Range("A6").Value = 3
$newSheet.Range("B6").Value = "FugaFuga"

# Style blank rows (2,5,7) with yellow fill
$newSheet.Range("A2:B2").Interior.Color = 65535
$newSheet.Range("A5:B5").Interior.Color = 65535
$newSheet.Range("A7:B7").Interior.Color = 65535

# Fix People sheet selection
$peopleSheet = $wb.Worksheets.Item("People")
$peopleSheet.Range("A1:D4").Select()

# Fix ResultCheckSheet selection / topLeftCell
$resultSheet = $wb.Worksheets.Item("ResultCheckSheet")
$resultSheet.Range("A1").Select()

$newSheet.Select()
